# Atlantic Coast Conf team order
#
# The "Atlantic Coast Conference" header/filler row had been left sitting at
# the top of the NCAA team list (row 4) instead of sorting into its correct
# alphabetical position among the team rows beneath it. Re-apply the sort on
# the data range (A4:G222, i.e. excluding the real column-header row 3) so it
# settles where it belongs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NCAA")

$dataRange = $ws.Range("A4:G222")
$keyRange = $ws.Range("A4:A222")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($keyRange)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Bring the NCAA sheet to the front (it was previously NBA) and leave the
# freshly-sorted range selected.
$ws.Activate()
$ws.Range("A3:G222").Select()
